$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove obsolete header columns (delete right-to-left so earlier
#     column letters stay valid while we work) ---

# "Phan loai" (U) and "HC Category" (V) -> deleted entirely
$ws.Range("U1:V1").EntireColumn.Delete()

# "Phut nghi phep","Phut tang ca 100%","Phut tang ca 150%","Phut tang ca dem","Phut nghi khong luong"
# (N:R) -> deleted; a single replacement column is inserted below
$ws.Range("N1:R1").EntireColumn.Delete()

# "So phut ca" (J) -> deleted entirely
$ws.Range("J1:J1").EntireColumn.Delete()

# "Cap bac" (G) -> deleted entirely
$ws.Range("G1:G1").EntireColumn.Delete()

# Insert the new replacement column ("Phut tang ca 200%") right before
# the surviving "Phut nghi khac" column (now column L)
$ws.Range("L1:L1").EntireColumn.Insert()
$ws.Range("L3").Value = "Phút tăng ca 200%"
$ws.Range("L3").Font.Bold = $true

# --- View: scroll right a bit and select the "Giờ vào" column ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$ws.Range("I1:I1048576").Select()
